# Apply review-feedback edits to the "HPCDATAMGM-1604" paragraph as tracked
# changes (insertions/deletions), matching the author of the original commit.

$word.UserName = "Menon, Sunita (NIH/NCI) [C]"
$d = $word.ActiveDocument
$d.TrackRevisions = $true

# ---------------------------------------------------------------------------
# Edit 1: " APIs to also provide" -> " APIs to provide"
#   (delete the word " also")
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(" also provide the paths of the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Edit 1 anchor text not found" }
$delRange1 = $d.Range($rng1.Start, $rng1.Start + 5)   # " also" -> 5 characters
$delRange1.Delete()

# ---------------------------------------------------------------------------
# Edit 2: "collections. Previously" -> "collections once they are completed. Previously"
#   (insert " once they are completed" right before the period)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("downloaded collections. Previously", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Edit 2 anchor text not found" }
$matchStart2 = $rng2.Start
$periodOffset2 = $rng2.Text.IndexOf(".")
$insPoint2 = $d.Range($matchStart2 + $periodOffset2, $matchStart2 + $periodOffset2)
$insPoint2.InsertBefore(" once they are completed")

# ---------------------------------------------------------------------------
# Edit 3: "these APIs provided only the paths" -> "these APIs only provided the paths"
#   (insert "only " before "provided"; delete " only" after "provided")
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("these APIs provided only the paths of the data objects", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Edit 3 anchor text not found" }
$matchStart3 = $rng3.Start
$matchText3 = $rng3.Text
$providedOffset3 = $matchText3.IndexOf("provided")
$onlyOffset3 = $matchText3.IndexOf("only") - 1   # include the leading space -> " only"

# delete " only" first so the earlier insert position (before "provided") is unaffected
$delRange3 = $d.Range($matchStart3 + $onlyOffset3, $matchStart3 + $onlyOffset3 + 5)
$delRange3.Delete()

$insPoint3 = $d.Range($matchStart3 + $providedOffset3, $matchStart3 + $providedOffset3)
$insPoint3.InsertBefore("only ")
